# Change the table style (Table Design gallery selection) on the three
# balance-sheet tables in this deck from the Google-Slides-imported
# custom style to the built-in PowerPoint table style
# "{4D7F8451-65D2-489C-B98B-010D39E13709}".
#
# Table.Style cannot be assigned directly (it's read-only); PowerPoint's
# object model requires Table.ApplyStyle("{GUID}") to change a table's
# style, exactly as picking a new style from the Table Design ribbon
# gallery would.

$p = $ppt.ActivePresentation

$oldStyleId = "{90D82F63-57D0-4330-ADD3-6C407542997A}"
$newStyleId = "{4D7F8451-65D2-489C-B98B-010D39E13709}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
